# Insert a new weekly price record for Jengibre (ginger) at row 46, pushing
# the existing rows 46..123 down to 47..124 (matches the author's diff, which
# is a straight row-insert: every old row N (46<=N<=123) reappears unchanged
# as new row N+1; the brand-new data lands in row 46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 46 downward by inserting a fresh row at position 46.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new observation.
$ws.Range("A46").Value = 6
$ws.Range("B46").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C46").Value = "Metropolitana"
$ws.Range("D46").Value = 45036
$ws.Range("E46").Value = 13
$ws.Range("F46").Value = 100114007
$ws.Range("G46").Value = "Jengibre"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 500
$ws.Range("K46").Value = 15000
$ws.Range("L46").Value = 16000
$ws.Range("M46").Value = 15540
$ws.Range("N46").Value = "`$/caja 13 kilos"
$ws.Range("O46").Value = "Perú"
$ws.Range("P46").Value = 1195
$ws.Range("Q46").Value = 13
$ws.Range("R46").Value = "Hortaliza"
